$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 266, shifting existing rows 266:277 down to 267:278.
$ws.Rows("266:266").Insert()

# Populate the newly inserted row 266 with the new weekly data point.
$ws.Range("A266").Value = 7
$ws.Range("B266").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C266").Value = "Ñuble"
$ws.Range("D266").Value = 44939
$ws.Range("E266").Value = 16
$ws.Range("F266").Value = 100112017
$ws.Range("G266").Value = "Apio"
$ws.Range("H266").Value = "Americana (o)"
$ws.Range("I266").Value = "Primera"
$ws.Range("J266").Value = 50
$ws.Range("K266").Value = 10000
$ws.Range("L266").Value = 10000
$ws.Range("M266").Value = 10000
$ws.Range("N266").Value = "$/docena de matas"
$ws.Range("O266").Value = "Provincia del Elquí"
$ws.Range("P266").Value = 1667
$ws.Range("Q266").Value = 6
$ws.Range("R266").Value = "Hortaliza"
